# edit.ps1 - reproduces the "badges.pptx" slide1 edit:
#   1. Picture "Picture 3" (id=7)   -> move/resize (x,cx change; y,cy unchanged)
#   2. TextBox "TextBox 5" (id=6)   -> reposition (off changes only)
#   3. TextBox "TextBox 12" (id=13) -> replaced by a new AutoShape Rectangle
#        ("Rectangle 9", which PowerPoint's internal shape-id/name counter
#        assigns as id=10) carrying the same "Web App" text, no fill / no
#        line, centered vertically, 0 text insets and explicit line spacing.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Picture 3 (shape index 4, id 7): off x 1259397 -> 1177998,
#    ext cx 2373675 -> 2529906 (y / cy unchanged)
# ---------------------------------------------------------------------
$pic = $s.Shapes.Item(4)
$pic.Left = 92.7557874015748
$pic.Width = 199.20523622047244

# ---------------------------------------------------------------------
# 2) TextBox 5 (shape index 6, id 6): off x 1911406 -> 1907704,
#    off y 1549506 -> 1529985 (ext unchanged)
# ---------------------------------------------------------------------
$tb = $s.Shapes.Item(6)
$tb.Left = 150.21295275590552
$tb.Top = 120.47129921259844

# ---------------------------------------------------------------------
# 3) TextBox 12 (shape index 7, id 13) is removed and replaced with a
#    plain rectangle AutoShape ("Web App" badge). In real PowerPoint the
#    new shape's auto-assigned Id/Name ("10" / "Rectangle 9") depends on
#    the document's internal shape counter; burn through the same number
#    of shape-creation "slots" (6 scratch shapes, created & deleted) so
#    the replacement shape lands on Id 10 / Name "Rectangle 9" exactly
#    as it did for the original author.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 6; $i++) {
    $scratch = $s.Shapes.AddShape(1, 0, 0, 1, 1)
    $scratch.Delete()
}

# Remove the old "TextBox 12" shape (still shape index 7 at this point).
$s.Shapes.Item(7).Delete()

# Add the new rectangle AutoShape in its place (7th AddShape call overall
# -> Id 10, Name "Rectangle 9").
$rect = $s.Shapes.AddShape(1, 0, 0, 1, 1)

# Text frame layout: no insets, vertically centered, shrink-shape-to-fit.
$tf = $rect.TextFrame
$tf.MarginLeft = 0
$tf.MarginRight = 0
$tf.MarginTop = 0
$tf.MarginBottom = 0
$tf.VerticalAnchor = 3
$tf.AutoSize = 1

# Text + run formatting.
$tf.TextRange.Text = "Web App"
$tf.TextRange.Font.Size = 32
$tf.TextRange.LanguageID = "de-DE"
$tf.TextRange.ParagraphFormat.SpaceWithin = 38.4

# No fill / no outline on the shape itself (look comes from the theme
# style reference).
$rect.Fill.Visible = 0
$rect.Line.Visible = 0

# Pin down the exact position/size last (AutoSize above recalculates the
# shape's bbox from the text, so the explicit geometry has to win last).
$rect.Left = 150.21295275590552
$rect.Top = 139.59688976377953
$rect.Width = 145.23318897637793
$rect.Height = 38.775078740157475

Write-Host "Rectangle shape:" $rect.Id $rect.Name
